# Actualización desde MV -datos-
# Appends 5 new daily rows (05-10-2021 .. 12-10-2021) to the bottom of the
# "Diaria" data table on the active sheet (rows 196-200, columns A:K).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @("05-10-2021", 36, 7, 8, -2, 0, -2, -7, -33, -14, 7),
    @("06-10-2021", 36, 7, 8, -2, 0, -2, -7, -33, -14, 7),
    @("07-10-2021", 36, 7, 8, -2, 0, -2, -7, -31, -16, 7),
    @("08-10-2021", 36, 7, 8, -2, 0, -2, -3, -30, -19, 7),
    @("12-10-2021", 36, 7, 8, -2, 0, -2, 2, -29, -25, 7)
)

$startRow = 196
$endRow = $startRow + $newRows.Count - 1

# Pre-format column A of the new rows as Text so Excel stores the
# dd-mm-yyyy labels as literal strings instead of auto-converting them
# into date serial numbers (matching the rest of the "Serie" column).
$ws.Range("A$startRow`:A$endRow").NumberFormat = "@"

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $rowNum = $startRow + $i
    $rowData = $newRows[$i]
    for ($col = 1; $col -le $rowData.Count; $col++) {
        $ws.Cells.Item($rowNum, $col).Value = $rowData[$col - 1]
    }
}

# Restore the default (unstyled) cell format so the new rows keep the
# same plain styling as the rest of the table.
$ws.Range("A$startRow`:K$endRow").Style = "Normal"
